# "Generate Report for handoff"
#
# The old source file "f3fe2d20-69bd-48b6-9903-50f9faa3a407.md" (previously
# "Ready for handoff") has been superseded in this report run by a new
# source file "0e6cf60e-1b1c-4937-ba19-c9d9fb74d796.md" whose handoff
# transform failed. ".localization-config" is still present (not to be
# localized) but now reported in row 3 instead of row 2. The de-de sheet's
# stale handoff-target (xlf) bookkeeping for the old file is reset back to
# defaults since the new file has not produced a handoff target yet.

$wb = $excel.ActiveWorkbook

$newFile = "0e6cf60e-1b1c-4937-ba19-c9d9fb74d796.md"
$oldFile = ".localization-config"
$commitHash = "dcbc587501b51dfdc95ab212bd266f485c76ca48"
$newFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitHash + "/e2e/" + $newFile
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/" + $commitHash + "/.localization-config"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

$wsOverview.Range("B3").Value = "Not to be localized"
$wsOverview.Range("C3").Value = "Not to be localized"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $newFileUrl, "", "", $newFile)
$wsOverview.Range("A2").Style = "Hyperlink"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", $oldFile)
$wsOverview.Range("A3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("B2").Value = "Handoff transform failed"
$wsZhCn.Range("D2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Ignored"

$wsZhCn.Range("B3").Value = "Not to be localized"
$wsZhCn.Range("D3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").Value = "Ignored"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $newFileUrl, "", "", $newFile)
$wsZhCn.Range("A2").Style = "Hyperlink"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configUrl, "", "", $oldFile)
$wsZhCn.Range("A3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("B2").Value = "Handoff transform failed"
$wsDeDe.Range("C2").Clear()
$wsDeDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Ignored"

$wsDeDe.Range("B3").Value = "Not to be localized"
$wsDeDe.Range("D3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").Value = "Ignored"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $newFileUrl, "", "", $newFile)
$wsDeDe.Range("A2").Style = "Hyperlink"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configUrl, "", "", $oldFile)
$wsDeDe.Range("A3").Style = "Hyperlink"
